$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '30.008.55'
$ws.Range('D3').Value = '1.883.52'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7362'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9994'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3159'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07163'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.67'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08312'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7551'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.401'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Value = '1.890.60'
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.59'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.146'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '29.970.99'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '248.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007832'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.150.96'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9972'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9994'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.880'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1568'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.262'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.67'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.043'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.472'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.555'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.532'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.179'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05316'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.80%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.246'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7669'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.89%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9974'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.720'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01955'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.755'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.51%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4560'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.031'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8783'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('D44').Value = '1.086.03'
$ws.Range('E44').Value = '  -1.17%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '72.29'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.851'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.531'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.532'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.029.75'
$ws.Range('E51').Value = '  -0.35%  '
